$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a pure numeric-looking string must be forced to
# Text format first, otherwise Excel auto-converts them to numbers (changing
# cell type from inline string to a number) instead of keeping them as text.
$textCells = @("D4","D5","D6","D8","D10","D14","D15","D20","D21","D23","D25","D26","D27","D32","D35","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '61.936.16'
$ws.Range('E2').Value = '  -2.21%  '
$ws.Range('D3').Value = '2.433.93'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '579.46'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('D6').Value = '142.67'
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.528'
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('D9').Value = '2.429.27'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -3.75%  '
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('E13').Value = '  -3.48%  '
$ws.Range('D14').Value = '26.28'
$ws.Range('E14').Value = '  -3.08%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0000171'
$ws.Range('E15').Value = '  -4.26%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.874.62'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '61.825.98'
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').Value = '2.423.91'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('E19').Value = '  -4.03%  '
$ws.Range('D20').Value = '7.05'
$ws.Range('E20').Value = '  -4.23%  '
$ws.Range('D21').Value = '329.03'
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('D23').Value = '1.95'
$ws.Range('E23').Value = '  -6.33%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '65.53'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '9.27'
$ws.Range('E26').Value = '  +4.84%  '
$ws.Range('D27').Value = '614.15'
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').Value = '2.554.49'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').Value = '0.0₃0945'
$ws.Range('E30').Value = '  -8.34%  '
$ws.Range('E31').Value = '  -6.28%  '
$ws.Range('D32').Value = '7.96'
$ws.Range('E32').Value = '  -3.94%  '
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('D35').Value = '4.88'
$ws.Range('E35').Value = '  -6.28%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  -7.23%  '
$ws.Range('D38').Value = '0.374'
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').Value = '151.06'
$ws.Range('E39').Value = '  +3.73%  '
$ws.Range('D40').Value = '18.28'
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('D41').Value = '5.21'
$ws.Range('E41').Value = '  -4.07%  '
$ws.Range('D42').Value = '1.75'
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('D43').Value = '42.45'
$ws.Range('E43').Value = '  +1.36%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '2.44'
$ws.Range('E45').Value = '  -9.84%  '
$ws.Range('D46').Value = '142.41'
$ws.Range('E46').Value = '  -4.48%  '
$ws.Range('D47').Value = '3.60'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('D48').Value = '0.0521'
$ws.Range('E48').Value = '  -2.93%  '
$ws.Range('D49').Value = '0.594'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').Value = '19.42'
$ws.Range('E50').Value = '  -8.76%  '
$ws.Range('D51').Value = '0.0902'
$ws.Range('E51').Value = '  -1.56%  '
